{"js": "// Update the date title and refresh the worksheet's division problems/answers.\n// The table keeps its original 20 rows x 5 columns shape: only the 5 rows\n// that actually contain text (rows 0, 4, 8, 12, 16) get new values; the\n// blank spacer rows in between are left untouched.\n\nconst body = context.document.body;\n\n// 1) Title paragraph: \"2024-07-04 Thursday\" -> \"2024-07-05 Friday\"\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  if (p.text && p.text.includes(\"2024-07-04 Thursday\")) {\n    p.insertText(\"2024-07-05 Friday\", \"Replace\");\n  }\n}\nawait context.sync();\n\n// 2) Table cell text updates, by (row index, col index) -> new text.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// New text values for the five populated rows, left-to-right.\nconst rowUpdates = {\n  0: [\"37\u00f76=6, 1\", \"24\u00f76=4, 0\", \"41\u00f72=20, 1\", \"42\u00f76=7, 0\", \"68\u00f75=13, 3\"],\n  4: [\"18\u00f76=3, 0\", \"56\u00f72=28, 0\", \"28\u00f77=4, 0\", \"14\u00f76=2, 2\", \"21\u00f78=2, 5\"],\n  8: [\"21\u00f75=4, 1\", \"36\u00f77=5, 1\", \"11\u00f73=3, 2\", \"33\u00f79=3, 6\", \"97\u00f72=48, 1\"],\n  12: [\"96\u00f78=12, 0\", \"88\u00f75=17, 3\", \"97\u00f74=24, 1\", \"82\u00f73=27, 1\", \"60\u00f72=30, 0\"],\n  16: [\"12\u00f78=1, 4\", \"82\u00f79=9, 1\", \"72\u00f79=8, 0\", \"41\u00f75=8, 1\", \"52\u00f79=5, 7\"],\n};\n\nfor (const rowIndexStr of Object.keys(rowUpdates)) {\n  const rowIndex = Number(rowIndexStr);\n  const newValues = rowUpdates[rowIndex];\n  const row = rows.items[rowIndex];\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  for (let c = 0; c < cells.items.length && c < newValues.length; c++) {\n    cells.items[c].value = newValues[c];\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date title and refresh the worksheet's division problems/answers.\n# The table keeps its original 20 rows x 5 columns shape: only the 5 rows\n# that actually contain text (rows 1, 5, 9, 13, 17 in 1-based COM indexing)\n# get new values; the blank spacer rows in between are left untouched.\n\n$d = $word.ActiveDocument\n\n# 1) Title paragraph: \"2024-07-04 Thursday\" -> \"2024-07-05 Friday\"\n$find = $d.Content.Find\n$find.Execute(\"2024-07-04 Thursday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2024-07-05 Friday\", 2)\n\n# 2) Table cell text updates, by (1-based row, 1-based col) -> new text.\n$t = $d.Tables.Item(1)\n\n$rowUpdates = @{\n    1  = @(\"37\u00f76=6, 1\", \"24\u00f76=4, 0\", \"41\u00f72=20, 1\", \"42\u00f76=7, 0\", \"68\u00f75=13, 3\")\n    5  = @(\"18\u00f76=3, 0\", \"56\u00f72=28, 0\", \"28\u00f77=4, 0\", \"14\u00f76=2, 2\", \"21\u00f78=2, 5\")\n    9  = @(\"21\u00f75=4, 1\", \"36\u00f77=5, 1\", \"11\u00f73=3, 2\", \"33\u00f79=3, 6\", \"97\u00f72=48, 1\")\n    13 = @(\"96\u00f78=12, 0\", \"88\u00f75=17, 3\", \"97\u00f74=24, 1\", \"82\u00f73=27, 1\", \"60\u00f72=30, 0\")\n    17 = @(\"12\u00f78=1, 4\", \"82\u00f79=9, 1\", \"72\u00f79=8, 0\", \"41\u00f75=8, 1\", \"52\u00f79=5, 7\")\n}\n\nforeach ($rowIndex in $rowUpdates.Keys) {\n    $newValues = $rowUpdates[$rowIndex]\n    $row = $t.Rows.Item($rowIndex)\n    for ($c = 1; $c -le $row.Cells.Count; $c++) {\n        $row.Cells.Item($c).Range.Text = $newValues[$c - 1]\n    }\n}\n"}
